$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark that currently sits in the middle
#    of the Portuguese introduction paragraph (splitting "...que prov" /
#    "avelmente...") and merge the two runs it separated back into a single
#    run, since the target XML has that text flowing through one run again.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$mergeFind = $d.Content
$mergeFind.Find.Execute("que provavelmente", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergeRange = $d.Range($mergeFind.Start, $mergeFind.End)
$mergeText = $mergeRange.Text
$mergeRange.Text = ""
$mergeRange.InsertAfter($mergeText)

# ---------------------------------------------------------------------------
# 2) Re-create the "_GoBack" bookmark spanning from the start of the Abstract
#    paragraph through "...Use this document as a ", splitting the run right
#    before "template" so that a <w:bookmarkStart> appears right after the
#    paragraph properties and a <w:bookmarkEnd> appears right before the
#    "template" run, exactly like Word leaves behind after the last edit.
# ---------------------------------------------------------------------------
$abstractFind = $d.Content
$abstractFind.Find.Execute("Abstract", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraStart = $abstractFind.Paragraphs(1).Range.Start

$templateFind = $d.Content
$templateFind.Find.Execute("Use this document as a ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$templateSplit = $templateFind.End

$goBackRange = $d.Range($paraStart, $templateSplit)
$d.Bookmarks.Add("_GoBack", $goBackRange)
